# hel-839: fix problème d'affichage mineurs dans excel
#
# The "Lisez-moi" sheet's C3 cell (indicateur "Taux de réalisation de
# l'activité") had its note text re-entered with a trailing line break,
# and several rows needed their heights bumped up so the (now slightly
# taller / re-wrapped) notes are fully visible again instead of being
# clipped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lisez-moi")

# --- C3: append a trailing newline to the note, preserving the existing
#     rich-text (bold label / plain body) runs. ---------------------------
$c3 = $ws.Cells.Item(3, 3)
$c3.Value2 = $c3.Value2 + "`n"

# Re-apply the bold runs (the four "Source : " / "Fréquence : " /
# "Mode de calcul :" / "Source(s) :" labels) that a plain Value2 write
# resets to the base (non-bold) font.
$c3.Characters(1, 9).Font.Bold = $true     # "Source : "
$c3.Characters(77, 12).Font.Bold = $true   # "Fréquence : "
$c3.Characters(99, 16).Font.Bold = $true   # "Mode de calcul :"
$c3.Characters(424, 11).Font.Bold = $true  # "Source(s) :"

# --- Row heights: re-fit the note rows (taller now / after re-wrap). ----
$ws.Rows.Item(2).RowHeight = 279
$ws.Rows.Item(3).RowHeight = 352.8
$ws.Rows.Item(4).RowHeight = 317.4
$ws.Rows.Item(5).RowHeight = 246
$ws.Rows.Item(8).RowHeight = 197.4
$ws.Rows.Item(13).RowHeight = 409.2
$ws.Rows.Item(14).RowHeight = 298.8
$ws.Rows.Item(15).RowHeight = 228
$ws.Rows.Item(16).RowHeight = 409.2
$ws.Rows.Item(17).RowHeight = 409.2
$ws.Rows.Item(18).RowHeight = 255
$ws.Rows.Item(19).RowHeight = 319.8
$ws.Rows.Item(20).RowHeight = 270.6
